$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new abandoned place entry (Hashima / Gunkanjima, Japan) ---
# Order matters: it reproduces the exact shared-string table layout that
# Excel produced (old slot 10 "Charleroi，Belgium" gets recycled for the
# first brand-new value typed in, "Hashima").
$ws.Range("B3").Value = "Hashima"
$ws.Range("C2").Value = "Charleroi,Belgium"
$ws.Range("C3").Value = "Nagasaki,Japan"
$ws.Range("E3").Value = "Hashima Island ( 端島 , or simply Hashima — -shima is a Japanese suffix for island ) , commonly called Gunkanjima (軍艦島; meaning Battleship Island ), is an abandoned island lying about 15 kilometers (9 miles) from the city of Nagasaki , in southern Japan. It is one of 505 uninhabited islands in Nagasaki Prefecture . The island's most notable features are its abandoned concrete buildings, undisturbed except by nature, and the surrounding sea wall ."
$ws.Range("D3").Value = "`nen.wikipedia.org"
$ws.Range("A3").Value = 2

# F3 should share the exact same "12/09/2019" text cell (stored as a shared
# string, not a date serial) that F2 already uses - copy it across instead
# of re-typing it so it keeps its plain (unstyled) number format.
$ws.Range("F2").Copy($ws.Range("F3"))

# Wrap text on the long-form columns of the new row (same visual style as
# the existing description cell, E2).
$ws.Range("D3").WrapText = $true
$ws.Range("E3").WrapText = $true

# Row heights to fit the wrapped content.
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 75

# Move / collapse the selection onto B3.
$ws.Range("B3").Select()

# Page setup: portrait, paper size 9 (A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
